$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "58.243.34"
$ws.Range("E2").Value2 = "  +3.11%  "

# Row 3
$ws.Range("D3").Value2 = "2.330.63"
$ws.Range("E3").Value2 = "  +0.18%  "

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value2 = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value2 = "  +0.05%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = "544.99"
$c.Style = "Normal"
$ws.Range("E5").Value2 = "  +6.36%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value2 = "134.95"
$c.Style = "Normal"
$ws.Range("E6").Value2 = "  +2.48%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value2 = "0.997"
$c.Style = "Normal"
$ws.Range("E7").Value2 = "  -0.18%  "

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value2 = "0.537"
$c.Style = "Normal"
$ws.Range("E8").Value2 = "  +0.70%  "

# Row 9
$ws.Range("D9").Value2 = "2.358.86"
$ws.Range("E9").Value2 = "  +1.25%  "

# Row 10
$ws.Range("E10").Value2 = "  +1.65%  "

# Row 11
$ws.Range("E11").Value2 = "  +1.26%  "

# Row 12
$ws.Range("E12").Value2 = "  +2.86%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value2 = "0.360"
$c.Style = "Normal"
$ws.Range("E13").Value2 = "  +6.82%  "

# Row 14
$ws.Range("B14").Value2 = "Avalanche"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value2 = "23.71"
$c.Style = "Normal"
$ws.Range("E14").Value2 = "  +0.88%  "

# Row 15
$ws.Range("B15").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value2 = "2.752.74"
$ws.Range("E15").Value2 = "  +0.41%  "

# Row 16
$ws.Range("D16").Value2 = "58.163.99"
$ws.Range("E16").Value2 = "  +3.06%  "

# Row 17
$ws.Range("E17").Value2 = "  +1.20%  "

# Row 18
$ws.Range("D18").Value2 = "2.343.90"
$ws.Range("E18").Value2 = "  +0.67%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value2 = "10.61"
$c.Style = "Normal"
$ws.Range("E19").Value2 = "  +1.68%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value2 = "335.65"
$c.Style = "Normal"
$ws.Range("E20").Value2 = "  +2.69%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value2 = "4.22"
$c.Style = "Normal"
$ws.Range("E21").Value2 = "  +2.16%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value2 = "6.74"
$c.Style = "Normal"
$ws.Range("E22").Value2 = "  -0.33%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value2 = "0.997"
$c.Style = "Normal"
$ws.Range("E23").Value2 = "  -0.23%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value2 = "62.02"
$c.Style = "Normal"
$ws.Range("E24").Value2 = "  +0.55%  "

# Row 25
$ws.Range("E25").Value2 = "  +4.35%  "

# Row 26
$ws.Range("B26").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C26").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value2 = "8.46"
$c.Style = "Normal"
$ws.Range("E26").Value2 = "  -2.13%  "

# Row 27
$ws.Range("B27").Value2 = "Binance-PegBSC-USD"
$ws.Range("C27").Value2 = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value2 = "0.996"
$c.Style = "Normal"
$ws.Range("E27").Value2 = "  -0.32%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value2 = "1.41"
$c.Style = "Normal"
$ws.Range("E28").Value2 = "  +7.43%  "

# Row 29
$ws.Range("E29").Value2 = "  +5.53%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value2 = "169.93"
$c.Style = "Normal"
$ws.Range("E30").Value2 = "  +1.53%  "

# Row 31
$ws.Range("D31").Value2 = "0.0₃0735"
$ws.Range("E31").Value2 = "  +2.37%  "

# Row 32
$ws.Range("E32").Value2 = "  +0.67%  "

# Row 33
$ws.Range("B33").Value2 = "SuiNetwork"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value2 = "1.03"
$c.Style = "Normal"
$ws.Range("E33").Value2 = "  +16.12%  "

# Row 34
$ws.Range("B34").Value2 = "EthereumClassic"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value2 = "18.52"
$c.Style = "Normal"
$ws.Range("E34").Value2 = "  +0.83%  "

# Row 35
$ws.Range("E35").Value2 = "  -0.02%  "

# Row 36
$ws.Range("B36").Value2 = "NEARProtocol"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value2 = "4.23"
$c.Style = "Normal"
$ws.Range("E36").Value2 = "  +8.15%  "

# Row 37
$ws.Range("B37").Value2 = "FirstDigitalUSD"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value2 = "0.997"
$c.Style = "Normal"
$ws.Range("E37").Value2 = "  -0.04%  "

# Row 38
$ws.Range("E38").Value2 = "  +0.37%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value2 = "1.64"
$c.Style = "Normal"
$ws.Range("E39").Value2 = "  +5.69%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value2 = "39.27"
$c.Style = "Normal"
$ws.Range("E40").Value2 = "  +2.25%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value2 = "150.62"
$c.Style = "Normal"
$ws.Range("E41").Value2 = "  +0.26%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value2 = "0.381"
$c.Style = "Normal"
$ws.Range("E42").Value2 = "  +2.25%  "

# Row 43
$ws.Range("B43").Value2 = "Bittensor"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value2 = "290.07"
$c.Style = "Normal"
$ws.Range("E43").Value2 = "  +3.98%  "

# Row 44
$ws.Range("B44").Value2 = "Filecoin"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value2 = "3.63"
$c.Style = "Normal"
$ws.Range("E44").Value2 = "  +1.93%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value2 = "19.30"
$c.Style = "Normal"
$ws.Range("E45").Value2 = "  +6.57%  "

# Row 46
$ws.Range("E46").Value2 = "  +0.08%  "

# Row 47
$ws.Range("E47").Value2 = "  +2.66%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value2 = "0.563"
$c.Style = "Normal"
$ws.Range("E48").Value2 = "  +1.00%  "

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value2 = "0.0218"
$c.Style = "Normal"
$ws.Range("E49").Value2 = "  +1.87%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = "17.57"
$c.Style = "Normal"
$ws.Range("E50").Value2 = "  +2.93%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value2 = "0.380"
$c.Style = "Normal"
$ws.Range("E51").Value2 = "  -0.02%  "
